# ActiTime Automation Final Level
# Rebuilds TC01 sheet content/layout (5 steps now include "Minimize Flyout
# Window" before Logout, plus 3 new columns: TestDataColumn /
# ExpectedResultColumn / ObjectMapColumn) and adds a new "testdata" sheet
# with Username/Password/URL sample data (URL cell carries a hyperlink).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TC01"

# ---------------------------------------------------------------------------
# TC01 : header row + 5 test-step rows, now with 3 extra columns
# ---------------------------------------------------------------------------

$ws1.Range("A1").Value = "TestScriptID"
$ws1.Range("B1").Value = "Description"
$ws1.Range("C1").Value = "MethodName"
$ws1.Range("D1").Value = "PackageClassName"

$ws1.Range("A2").Value = "TS01"
$ws1.Range("B2").Value = "Navigate the Application URL"
$ws1.Range("C2").Value = "navigate"
$ws1.Range("D2").Value = "com.sgtesting.actitime.tests.Initialize"

$ws1.Range("A3").Value = "TS02"
$ws1.Range("B3").Value = "Login into the Application"
$ws1.Range("C3").Value = "login"
$ws1.Range("D3").Value = "com.sgtesting.actitime.tests.LoginLogout"

$ws1.Range("A4").Value = "TS03"
$ws1.Range("B4").Value = "Minimize Flyout Window"
$ws1.Range("C4").Value = "minimizeFlyOutWindow"
$ws1.Range("D4").Value = "com.sgtesting.actitime.tests.HomePage"

$ws1.Range("A5").Value = "TS04"
$ws1.Range("B5").Value = "Logout from the Application"
$ws1.Range("C5").Value = "logout"
$ws1.Range("D5").Value = "com.sgtesting.actitime.tests.LoginLogout"

$ws1.Range("A6").Value = "TS05"
$ws1.Range("B6").Value = "Close Application"
$ws1.Range("C6").Value = "closeApplication"
$ws1.Range("D6").Value = "com.sgtesting.actitime.tests.Initialize"

# Extra columns: TestDataColumn / ExpectedResultColumn / ObjectMapColumn
$ws1.Range("E1").Value = "TestDataColumn"
$ws1.Range("F1").Value = "ExpectedResultColumn"
$ws1.Range("G1").Value = "ObjectMapColumn"

$ws1.Range("E2").Value = "URL"
$ws1.Range("E3").Value = "Username,Password"

$ws1.Range("F2").Value = "logintitle"
$ws1.Range("F3").Value = "homepagetext"
$ws1.Range("F5").Value = "logintitle"

$ws1.Range("G3").Value = "txtloginusername,txtloginpassword,btnloginbutton"
$ws1.Range("G5").Value = "lnkhomelogout"

# copy the bold/boxed header style (A1) onto the new header cells (values
# were already written above; PasteSpecial(formats) leaves them intact)
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("E1:G1").PasteSpecial(-4122) | Out-Null

# copy the boxed body style (A2) onto the new E/F/G cells so the borders /
# font match the rest of the table
$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("E2:E6").PasteSpecial(-4122) | Out-Null
$ws1.Range("F2:F6").PasteSpecial(-4122) | Out-Null
$ws1.Range("G2:G6").PasteSpecial(-4122) | Out-Null

$ws1.Range("G3:G4").WrapText = $true
$ws1.Rows.Item(3).RowHeight = 29

# column widths (closest achievable match to the target bestFit widths -
# this runtime quantizes ColumnWidth to 1/6 character-width steps)
$ws1.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws1.Columns.Item(2).ColumnWidth = 26.333333333333336
$ws1.Columns.Item(3).ColumnWidth = 20.0
$ws1.Columns.Item(4).ColumnWidth = 34.33333333333333
$ws1.Columns.Item(5).ColumnWidth = 17.0
$ws1.Columns.Item(6).ColumnWidth = 20.0
$ws1.Columns.Item(7).ColumnWidth = 32.666666666666664

$ws1.Range("A1").Select()
$excel.ActiveWindow.Zoom = 130

# ---------------------------------------------------------------------------
# testdata sheet : Username / Password / URL sample row, URL is a hyperlink
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "testdata"

$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "URL"
$ws2.Range("A1:C1").Font.Bold = $true

$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "manager"
$ws2.Range("C2").Value = "http://localhost/login.do"

$ws2.Hyperlinks.Add($ws2.Range("C2"), "http://localhost/login.do") | Out-Null

$ws2.Columns.Item(1).ColumnWidth = 8.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 8.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 25.333333333333336

$ws2.Range("A2:C2").Select()
$excel.ActiveWindow.Zoom = 160

$ws1.Select()
